# AAP-36 + AAP-34: Import and export
# Split the single "Service Coverage" column into two columns:
#   "Service Coverage Region"  (existing column G, text updated)
#   "Service Coverage Country" (brand new column, inserted at H)
# Everything that used to live from column H onward (Global Focal Point, ...,
# Inter-Agency CFM Resources) shifts one column to the right automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before H - this shifts H:R -> I:S and extends
# the formatted/width column range + used dimension automatically.
$ws.Columns("H").Insert()

# Rename the existing "Service Coverage" header (column G) and give the new
# column (H) its own header text.
$ws.Range("G1").Value2 = "Service Coverage Region"
$ws.Range("H1").Value2 = "Service Coverage Country"

# The _FilterDatabase defined name covered A1:M1 before; it now needs to
# cover one extra column (A1:N1) to match the widened header row.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Services!`$A`$1:`$N`$1"
    }
}

# Leave the cursor on the newly added header cell, matching the saved
# selection state of the edited workbook.
$ws.Range("H1").Select() | Out-Null
